$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend header row 1 with two new columns (P=14, Q=15), matching the
# existing bold/bordered header style (copied from O1).
$ws.Range("O1").Copy() | Out-Null
$ws.Range("P1:Q1").PasteSpecial(-4122) | Out-Null
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

# Updated simulation results for columns B:I, rows 2-25.
$ws.Range("B2").Value = 3.457262170408171
$ws.Range("C2").Value = 1.007702210158641
$ws.Range("D2").Value = 0.3492156199738474
$ws.Range("E2").Value = 1.343180230842464
$ws.Range("F2").Value = 6.463794665469322
$ws.Range("G2").Value = 0.0007918837209675011
$ws.Range("H2").Value = 0.009995295743722019
$ws.Range("I2").Value = 0.002670933310996659

$ws.Range("B3").Value = 2.997853015548628
$ws.Range("C3").Value = 0.869853693225707
$ws.Range("D3").Value = 0.3064966316520383
$ws.Range("E3").Value = 1.155751310872887
$ws.Range("F3").Value = 5.654258087503649
$ws.Range("G3").Value = 0.0007995267261627653
$ws.Range("H3").Value = 0.00566132105441941
$ws.Range("I3").Value = 0.0006838291904367466

$ws.Range("B4").Value = 2.717878933285135
$ws.Range("C4").Value = 0.7870470786523356
$ws.Range("D4").Value = 0.2802437474734631
$ws.Range("E4").Value = 1.042418103720038
$ws.Range("F4").Value = 5.159655295851422
$ws.Range("G4").Value = 0.00080432863316101
$ws.Range("H4").Value = 0.003562661497822361
$ws.Range("I4").Value = 0.0004342247631967666

$ws.Range("B5").Value = 2.604150420224698
$ws.Range("C5").Value = 0.7552167163809997
$ws.Range("D5").Value = 0.2683125171909211
$ws.Range("E5").Value = 0.9965441341777392
$ws.Range("F5").Value = 4.944718871251297
$ws.Range("G5").Value = 0.0008063365970838446
$ws.Range("H5").Value = 0.00283109922356406
$ws.Range("I5").Value = 0.000622890166188661

$ws.Range("B6").Value = 2.585207879043878
$ws.Range("C6").Value = 0.7518079961090791
$ws.Range("D6").Value = 0.2648516623235935
$ws.Range("E6").Value = 0.9888719871438241
$ws.Range("F6").Value = 4.892255946361189
$ws.Range("G6").Value = 0.0008066988388415279
$ws.Range("H6").Value = 0.002711539273524188
$ws.Range("I6").Value = 0.0007586302123687716

$ws.Range("B7").Value = 2.716100113317282
$ws.Range("C7").Value = 0.7916813407236987
$ws.Range("D7").Value = 0.2760110113271139
$ws.Range("E7").Value = 1.041580167724106
$ws.Range("F7").Value = 5.110605349627804
$ws.Range("G7").Value = 0.0008044293528601984
$ws.Range("H7").Value = 0.00353482695453966
$ws.Range("I7").Value = 0.0006621355425160402

$ws.Range("B8").Value = 3.297930760598206
$ws.Range("C8").Value = 0.9664754913100921
$ws.Range("D8").Value = 0.3289524093292329
$ws.Range("E8").Value = 1.277802295241244
$ws.Range("F8").Value = 6.121670233424283
$ws.Range("G8").Value = 0.0007945947571169861
$ws.Range("H8").Value = 0.00834023203055928
$ws.Range("I8").Value = 0.001998973943032389

$ws.Range("B9").Value = 4.461276162754871
$ws.Range("C9").Value = 1.317628976875312
$ws.Range("D9").Value = 0.4409129175473083
$ws.Range("E9").Value = 1.758791490145597
$ws.Range("F9").Value = 8.221729476999883
$ws.Range("G9").Value = 0.0007759650060087116
$ws.Range("H9").Value = 0.02293966575955908
$ws.Range("I9").Value = 0.01318788530029646

$ws.Range("B10").Value = 5.332556541617578
$ws.Range("C10").Value = 1.589096750276838
$ws.Range("D10").Value = 0.5003877811189454
$ws.Range("E10").Value = 2.021280822810638
$ws.Range("F10").Value = 9.52940342855976
$ws.Range("G10").Value = 0.0007632755018147897
$ws.Range("H10").Value = 0.03654587955837574
$ws.Range("I10").Value = 0.02786243332301552

$ws.Range("B11").Value = 5.696287717332439
$ws.Range("C11").Value = 1.68929629314573
$ws.Range("D11").Value = 0.3475534396946358
$ws.Range("E11").Value = 1.319623735495057
$ws.Range("F11").Value = 7.984121445527933
$ws.Range("G11").Value = 0.0007625219121048961
$ws.Range("H11").Value = 0.05102710592101189
$ws.Range("I11").Value = 0.02986244529849014

$ws.Range("B12").Value = 5.821453138132654
$ws.Range("C12").Value = 1.71081905899149
$ws.Range("D12").Value = 0.2378295168972926
$ws.Range("E12").Value = 0.8095450099342685
$ws.Range("F12").Value = 6.621540193137264
$ws.Range("G12").Value = 0.0007639042150792652
$ws.Range("H12").Value = 0.08537059510579326
$ws.Range("I12").Value = 0.02834500312373667

$ws.Range("B13").Value = 5.769062860697261
$ws.Range("C13").Value = 1.681431827745541
$ws.Range("D13").Value = 0.1489218976295632
$ws.Range("E13").Value = 0.4192333535295205
$ws.Range("F13").Value = 5.257492853800528
$ws.Range("G13").Value = 0.0007670939994285801
$ws.Range("H13").Value = 0.1363447470369579
$ws.Range("I13").Value = 0.02441139653021018

$ws.Range("B14").Value = 5.653258874752623
$ws.Range("C14").Value = 1.640501665964337
$ws.Range("D14").Value = 0.09914364752189897
$ws.Range("E14").Value = 0.2200208786317432
$ws.Range("F14").Value = 4.320288365004302
$ws.Range("G14").Value = 0.0007700778533965794
$ws.Range("H14").Value = 0.18241339872894
$ws.Range("I14").Value = 0.02096641771609153

$ws.Range("B15").Value = 5.587027398541977
$ws.Range("C15").Value = 1.621696695691526
$ws.Range("D15").Value = 0.08790478637317278
$ws.Range("E15").Value = 0.180148105521539
$ws.Range("F15").Value = 4.065077586257132
$ws.Range("G15").Value = 0.0007712711082461959
$ws.Range("H15").Value = 0.1938160907104987
$ws.Range("I15").Value = 0.01972158660848322

$ws.Range("B16").Value = 5.224085017559617
$ws.Range("C16").Value = 1.516492032032318
$ws.Range("D16").Value = 0.08526566423042681
$ws.Range("E16").Value = 0.1700411051112454
$ws.Range("F16").Value = 3.828971836814645
$ws.Range("G16").Value = 0.00077591099244641
$ws.Range("H16").Value = 0.1772398316277872
$ws.Range("I16").Value = 0.01491254195543057

$ws.Range("B17").Value = 5.007352928785792
$ws.Range("C17").Value = 1.457340512143503
$ws.Range("D17").Value = 0.1063005318635462
$ws.Range("E17").Value = 0.2527426306022065
$ws.Range("F17").Value = 4.119620197415884
$ws.Range("G17").Value = 0.0007779609381707517
$ws.Range("H17").Value = 0.1380845818779335
$ws.Range("I17").Value = 0.01298819220147251

$ws.Range("B18").Value = 4.891181072940014
$ws.Range("C18").Value = 1.426784906721878
$ws.Range("D18").Value = 0.1630917184647416
$ws.Range("E18").Value = 0.488672800946965
$ws.Range("F18").Value = 4.989782952390868
$ws.Range("G18").Value = 0.0007777175351106891
$ws.Range("H18").Value = 0.08624958060428156
$ws.Range("I18").Value = 0.0128610273449894

$ws.Range("B19").Value = 4.865410637195282
$ws.Range("C19").Value = 1.433989826580387
$ws.Range("D19").Value = 0.2585924094822332
$ws.Range("E19").Value = 0.9316599478342198
$ws.Range("F19").Value = 6.311405617291229
$ws.Range("G19").Value = 0.0007753471943131185
$ws.Range("H19").Value = 0.04486607482094485
$ws.Range("I19").Value = 0.01491509726212481

$ws.Range("B20").Value = 5.100185226391659
$ws.Range("C20").Value = 1.532472938516833
$ws.Range("D20").Value = 0.4702190910753643
$ws.Range("E20").Value = 1.945769640609583
$ws.Range("F20").Value = 9.0203272447026
$ws.Range("G20").Value = 0.0007668416119520997
$ws.Range("H20").Value = 0.03244554102579755
$ws.Range("I20").Value = 0.02379408360536761

$ws.Range("B21").Value = 5.784892976581773
$ws.Range("C21").Value = 1.750257957264921
$ws.Range("D21").Value = 0.5534414014092022
$ws.Range("E21").Value = 2.318371252717682
$ws.Range("F21").Value = 10.45960704905036
$ws.Range("G21").Value = 0.000756201365715177
$ws.Range("H21").Value = 0.0464444912688986
$ws.Range("I21").Value = 0.03888522512563775

$ws.Range("B22").Value = 6.238143054835405
$ws.Range("C22").Value = 1.888477988335353
$ws.Range("D22").Value = 0.6034973891750326
$ws.Range("E22").Value = 2.513470465521181
$ws.Range("F22").Value = 11.35561038052936
$ws.Range("G22").Value = 0.0007494717359220136
$ws.Range("H22").Value = 0.05622256578920748
$ws.Range("I22").Value = 0.05032041316143676

$ws.Range("B23").Value = 5.996341369730999
$ws.Range("C23").Value = 1.808243448853432
$ws.Range("D23").Value = 0.5820632899846885
$ws.Range("E23").Value = 2.409407682906746
$ws.Range("F23").Value = 10.93547678253219
$ws.Range("G23").Value = 0.0007529736297253577
$ws.Range("H23").Value = 0.05097258498324031
$ws.Range("I23").Value = 0.04393539907997379

$ws.Range("B24").Value = 5.093396454990398
$ws.Range("C24").Value = 1.522929633599688
$ws.Range("D24").Value = 0.494095960497134
$ws.Range("E24").Value = 2.024357664267967
$ws.Range("F24").Value = 9.282577626360364
$ws.Range("G24").Value = 0.0007664096614360714
$ws.Range("H24").Value = 0.03328234609949288
$ws.Range("I24").Value = 0.02388743652415215

$ws.Range("B25").Value = 4.142656613778001
$ws.Range("C25").Value = 1.229542174013773
$ws.Range("D25").Value = 0.4026513398619613
$ws.Range("E25").Value = 1.625711697670212
$ws.Range("F25").Value = 7.560761790821061
$ws.Range("G25").Value = 0.0007810655889403182
$ws.Range("H25").Value = 0.01826779157226266
$ws.Range("I25").Value = 0.009293847328576632

# New trailing columns P:Q, rows 2-25, filled with 0.
$ws.Range("P2:Q25").Value = 0
